$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("depo")

# --- Update depo name text (column C) to shortened / renamed values ---
$ws.Range("C1").Value = "Depo Lokomotif THB"
$ws.Range("C2").Value = "Depo Lokomotif CPN"
$ws.Range("C3").Value = "Depo Kereta JAKK"
$ws.Range("C4").Value = "Depo Gerbong JAKG"
$ws.Range("C5").Value = "Depo Kereta RK"
$ws.Range("C6").Value = "Depo Kereta JAKRI"
$ws.Range("C9").Value = "Depo Lokomotif BD"
$ws.Range("C10").Value = "Depo Kereta BD"
$ws.Range("C11").Value = "Depo Lokomotif CN"
$ws.Range("C12").Value = "Depo Kereta CN"
$ws.Range("C13").Value = "Depo Lokomotif SMC"
$ws.Range("C14").Value = "Depo Kereta SMC"
$ws.Range("C17").Value = "Depo Lokomotif PWT"
$ws.Range("C18").Value = "Depo Kereta PWT"
$ws.Range("C19").Value = "Depo Kereta KTA"
$ws.Range("C20").Value = "Depo Gerbong MA"
$ws.Range("C23").Value = "Depo Lokomotif YK"
$ws.Range("C24").Value = "Depo Kereta YK"
$ws.Range("C25").Value = "Depo Lokomotif SLO"
$ws.Range("C26").Value = "Depo Kereta SLO"
$ws.Range("C28").Value = "Depo Lokomotif MN"
$ws.Range("C29").Value = "Depo Kereta BL"
$ws.Range("C30").Value = "Depo Mekanik MN"
$ws.Range("C31").Value = "Depo Lokomotif SDT"
$ws.Range("C32").Value = "Depo Kereta SDT"
$ws.Range("C33").Value = "Depo Kereta SBI"
$ws.Range("C35").Value = "Depo Lokomotif ML"
$ws.Range("C36").Value = "Depo Kereta ML"
$ws.Range("C37").Value = "Depo Gerbong SDT"
$ws.Range("C38").Value = "Depo Lokomotif JR"
$ws.Range("C39").Value = "Depo Kereta BW"
$ws.Range("C40").Value = "Depo Lokomotif MDN"
$ws.Range("C41").Value = "Depo Kereta MDN"
$ws.Range("C42").Value = "Depo KRD MDN"
$ws.Range("C44").Value = "Depo Lokomotif PD"
$ws.Range("C45").Value = "Depo Kereta PD "
$ws.Range("C46").Value = "Depo Lokomotif KPT"
$ws.Range("C47").Value = "Depo Kereta KPT"
$ws.Range("C48").Value = "Depo Gerbong KPT"
$ws.Range("C51").Value = "DEPO LRT JAKABARING"
$ws.Range("C52").Value = "Depo KRL BukitDuri"
$ws.Range("C56").Value = "Depo Lokomotif THN"
$ws.Range("C57").Value = "Depo Kereta TNK"
$ws.Range("C58").Value = "Depo Lokomotif RJS"
$ws.Range("C59").Value = "Depo Kereta RJS"
$ws.Range("C60").Value = "Depo Gerbong RJS"
$ws.Range("C61").Value = "Depo Lokomotif TNK"
$ws.Range("C62").Value = "Depo Gerbong THN"
$ws.Range("C63").Value = "Depo MRT LebakBulus"
$ws.Range("C64").Value = "DEPO LRT KELAPAGADING"
$ws.Range("C65").Value = "Depo APMS BandaraSoetta"

# --- Append 5 newly added depo rows at the bottom of the depo sheet ---
$ws.Range("A66").Value = 79
$ws.Range("B66").Value = "DAOP 3 CIREBON"
$ws.Range("C66").Value = "Depo Gerbong AWN"
$ws.Range("D66").Value = "PT Kereta Api Indonesia (Persero)"
$ws.Range("E66").Value = "Depo Sarana Tanpa Penggerak"
$ws.Range("A67").Value = 80
$ws.Range("B67").Value = "DAOP 6 YOGYAKARTA"
$ws.Range("C67").Value = "DEPO GERBONG RWL"
$ws.Range("D67").Value = "PT Kereta Api Indonesia (Persero)"
$ws.Range("E67").Value = "Depo Sarana Tanpa Penggerak"
$ws.Range("A68").Value = 81
$ws.Range("B68").Value = "DIVRE II SUMATERA BARAT"
$ws.Range("C68").Value = "DEPO GERBONG BKP"
$ws.Range("D68").Value = "PT Kereta Api Indonesia (Persero)"
$ws.Range("E68").Value = "Depo Sarana Tanpa Penggerak"
$ws.Range("A69").Value = 82
$ws.Range("B69").Value = "DIVRE III SUMATERA SELATAN"
$ws.Range("C69").Value = "Depo Gerbong SIG"
$ws.Range("D69").Value = "PT Kereta Api Indonesia (Persero)"
$ws.Range("E69").Value = "Depo Sarana Tanpa Penggerak"
$ws.Range("A70").Value = 83
$ws.Range("B70").Value = "DIVRE III SUMATERA SELATAN"
$ws.Range("C70").Value = "DEPO GERBONG MRL"
$ws.Range("D70").Value = "PT Kereta Api Indonesia (Persero)"
$ws.Range("E70").Value = "Depo Sarana Tanpa Penggerak"

# --- Add a new staging worksheet "Sheet1" listing the newly added depos ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet1"

$ws2.Range("B2").Value = "tambah depo baru"
$ws2.Range("A4").Value = "Daop 3"
$ws2.Range("B4").Value = "Depo Gerbong AWN"
$ws2.Range("A5").Value = "DAOP 6"
$ws2.Range("B5").Value = "DEPO GERBONG RWL"
$ws2.Range("A6").Value = "Divre 2"
$ws2.Range("B6").Value = "DEPO GERBONG BKP"
$ws2.Range("A7").Value = "Divre 3"
$ws2.Range("B7").Value = "Depo Gerbong Simpang"
$ws2.Range("A8").Value = "Divre 3"
$ws2.Range("B8").Value = "DEPO GERBONG MRL"
